# Update "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 203ff4b0-... row on both language sheets, as part
# of regenerating the Handback status report.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-20 06:34:25"
$wsZh.Range("H3").Value = "2016-03-20 06:34:45"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-20 06:34:28"
$wsDe.Range("H3").Value = "2016-03-20 06:34:51"
